$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix the "Objetivos:" value (row 10) - it currently (wrongly) holds the
#    "Docentes responsaveis" text; replace it with the real objectives text.
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = 'Aprofundar os conceitos técnicos fundamentais de um curso de Engenharia de Produção, tendo em vista a sua formação generalista voltada para os mais diversos tipos de sistemas de produção.'
$ws.Range("C10").Value = 'Aprofundar os conceitos técnicos fundamentais de um curso de Engenharia de Produção, tendo em vista a sua formação generalista voltada para os mais diversos tipos de sistemas de produção.'

# ---------------------------------------------------------------------------
# 2) Insert a new row right after row 12 ("Docentes responsaveis:") that
#    carries the value that used to incorrectly live on the "Objetivos:" row.
#    This shifts every following row down by one (old row 13 -> new row 14,
#    old row 21 -> new row 22, etc.)
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Insert()

# The inserted row copies formatting from the row above, which leaves a
# stray styled-but-empty A13 cell; drop it since this row has no label.
$ws.Range("A13").Clear()

# Pull the B/C formatting from the row below (same column convention) so the
# new cells land on the existing "column" styles instead of synthesizing new
# ones from individual Font/WrapText/Alignment property writes.
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Value = '5840535 - Messias Borges Silva'

$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C13").Value = '5840535 - Messias Borges Silva'

# ---------------------------------------------------------------------------
# 3) Correct the values that were each shifted one slot away from the label
#    they actually belong to (rows renumbered after the insert above).
# ---------------------------------------------------------------------------

# Programa resumido: (was "Semestral")
$ws.Range("B14").Value = '1 – Planejamento e Controle da Qualidade2 – Melhoramentos da Produção3 – Desafios da produção4 – Controle da Qualidade'
$ws.Range("C14").Value = '1 – Planejamento e Controle da Qualidade2 – Melhoramentos da Produção3 – Desafios da produção4 – Controle da Qualidade'

# Programa: (was "01/01/2018")
$ws.Range("B16").Value = '1 – Planejamento e Controle da QualidadeIntrodução. Planejamento e Controle da qualidade.2 – Melhoramentos da ProduçãoIntrodução. Medidas e melhoramentos de desempenho. Prevenção e Recuperação de falhas. Administração da Qualidade Total.3 – Desafios da produçãoIntrodução. Tipo e formas de estratégias.4 - CONTROLE DA QUALIDADEAs Sete Ferramentas da Qualidade: Diagrama de Ishikawa, Histograma, Folha de Verificação, Estratificação, Diagrama de Pareto, Diagrama de Dispersão, Gráficos de Controle. Círculos de Controle da Qualidade'
$ws.Range("C16").Value = '1 – Planejamento e Controle da QualidadeIntrodução. Planejamento e Controle da qualidade.2 – Melhoramentos da ProduçãoIntrodução. Medidas e melhoramentos de desempenho. Prevenção e Recuperação de falhas. Administração da Qualidade Total.3 – Desafios da produçãoIntrodução. Tipo e formas de estratégias.4 - CONTROLE DA QUALIDADEAs Sete Ferramentas da Qualidade: Diagrama de Ishikawa, Histograma, Folha de Verificação, Estratificação, Diagrama de Pareto, Diagrama de Dispersão, Gráficos de Controle. Círculos de Controle da Qualidade'

# Método: (was "5840535 - Messias Borges Silva")
$ws.Range("B19").Value = 'Aulas Expositivas; trabalhos e seminários.'
$ws.Range("C19").Value = 'Aulas Expositivas; trabalhos e seminários.'

# Critério: (was "Aulas Expositivas; trabalhos e seminários.")
$ws.Range("B20").Value = 'MF = (0,30*P1 + 0,30*P2 + 0,40*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários.'
$ws.Range("C20").Value = 'MF = (0,30*P1 + 0,30*P2 + 0,40*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários.'

# Norma de recuperação: (was "MF = (0,30*P1 ...")
$ws.Range("B21").Value = 'NF = (MF + PR)/2, onde PR é uma prova de recuperação.'
$ws.Range("C21").Value = 'NF = (MF + PR)/2, onde PR é uma prova de recuperação.'

# Bibliografia: (was "NF = (MF + PR)/2 ...")
$ws.Range("B22").Value = 'SLACK, N. et al. Administração da produção. São Paulo: Atlas, 2002. VENANZI, D; SILVA, O.R., Gerenciamento da Produçao e Operaçoes, LTC, 2014Textos complementares serão usados durante o curso.'
$ws.Range("C22").Value = 'SLACK, N. et al. Administração da produção. São Paulo: Atlas, 2002. VENANZI, D; SILVA, O.R., Gerenciamento da Produçao e Operaçoes, LTC, 2014Textos complementares serão usados durante o curso.'
